$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tanque")

$ws.Range("F2").Value = "Divergência entre o SPED(4783,00) e o relatório(9789,90)!"
$ws.Range("F3").Value = "Divergência entre o SPED(4783,00) e o relatório(3243,18)!"
$ws.Range("F4").Value = "Divergência entre o SPED(4783,00) e o relatório(4535,18)!"
$ws.Range("F5").Value = "Divergência entre o SPED(4783,00) e o relatório(1627,98)!"
$ws.Range("F6").Value = "Divergência entre o SPED(4783,00) e o relatório(2767,01)!"
